$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column, matching the formatting of the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add "Save" values for rows 2-7
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
